# TC01-Create Crisis.xlsx -- test case restructuring by Gilana and Waseem
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Options"
# ---------------------------------------------------------------
$opt = $wb.Worksheets.Item("Options")

# Row 2: fix duplicated "Starts with number" label in G2, clear H2
$opt.Range("G2").Value = "Contains special character"
$opt.Range("H2").Value = ""
# Row 2 no longer needs the taller (wrapped) row height
$opt.Rows.Item(2).EntireRow.AutoFit()

# Row 3: Crisis Type example values replaced with Selected / Not selected
$opt.Range("B3").Value = "Selected"
$opt.Range("C3").Value = "<Not selected>"

# Row 4: new "Crisis Location Name" option row
$opt.Range("A4").Value = "Crisis Location Name"
$opt.Range("B4").Value = "Space"
$opt.Range("C4").Value = "Letters"
$opt.Range("D4").Value = "Numbers"
$opt.Range("E4").Value = "Letters and Numbers"
$opt.Rows.Item(4).RowHeight = 30

# Column widths for the newly used E/F/G columns
$opt.Columns.Item(5).ColumnWidth = 18.165
$opt.Columns.Item(6).ColumnWidth = 19.165
$opt.Columns.Item(7).ColumnWidth = 24.335

# ---------------------------------------------------------------
# Sheet "Test Cases"
# ---------------------------------------------------------------
$tc = $wb.Worksheets.Item("Test Cases")

# Row 1 explanation text (shared strings got reshuffled, text unchanged)
$tc.Range("A1").Value = "Explanation:"
$tc.Range("B1").Value = "Testcase ID: First two numbers is the use case number, second two numbers are remaining test case number. Input Explanation: O1 means option 1 will be selected for this case. Ore is the example to use in test case."

# Row 3: new "Crisis Location Name" column header
$tc.Range("E3").Value = "Crisis Location Name"

# Row 4 (TC01.01): location name column added, input changed
$tc.Range("C4").Value = "O2:Abc"
$tc.Range("D4").Value = "O2:Fire"
$tc.Range("E4").Value = "O3:Abc 23"
$tc.Rows.Item(4).RowHeight = 33.75

# Row 5 (TC01.02): location name column added
$tc.Range("E5").Value = "O3:Abcdfg"
$tc.Rows.Item(5).RowHeight = 21

# Row 6 (TC01.03): new test case - crisis without location name
$tc.Range("A6").Value = "TC01.03"
$tc.Range("B6").Value = "Create a crisis with out location name"
$tc.Range("C6").Value = "O1:Abc"
$tc.Range("D6").Value = "O2:Fire"
$tc.Range("E6").Value = "O1:<space>"
$tc.Range("K6").Value = "Error Message"
$tc.Rows.Item(6).RowHeight = 24.75

# Row 7 (TC01.04): new test case - crisis with only numbers in location name
$tc.Range("A7").Value = "TC01.04"
$tc.Range("B7").Value = "Create a crisis with only numbers in location name"
$tc.Range("C7").Value = "O2:Ab"
$tc.Range("D7").Value = "O1:Earthquake"
$tc.Range("E7").Value = "O3:3278648"
$tc.Range("K7").Value = "Error Message"
$tc.Rows.Item(7).RowHeight = 26.25

# Row 8 (TC01.05): new test case - crisis without selecting Crisis Type
$tc.Range("A8").Value = "TC01.05"
$tc.Range("B8").Value = "Create crisis without selecting Crisis Type"
$tc.Range("C8").Value = "O2:Abc"
$tc.Range("D8").Value = "O2:<not selected>"
$tc.Range("E8").Value = "O3:Abc 23"
$tc.Range("K8").Value = "Error Message"

# Column widths across the newly populated columns
$tc.Columns.Item(2).ColumnWidth = 45.5
$tc.Columns.Item(4).ColumnWidth = 16.5
$tc.Columns.Item(5).ColumnWidth = 20.0
$tc.Columns.Item(6).ColumnWidth = 9.17
$tc.Columns.Item(7).ColumnWidth = 9.67
$tc.Columns.Item(8).ColumnWidth = 7.67
$tc.Columns.Item(9).ColumnWidth = 6.67
$tc.Columns.Item(10).ColumnWidth = 7.67

# ---------------------------------------------------------------
# Sheet selections / active tab
# ---------------------------------------------------------------
# Options: select C4, no longer the active tab
$opt.Range("C4").Select()

# Sheet3: select C27
$sheet3 = $wb.Worksheets.Item("Sheet3")
$sheet3.Range("C27").Select()

# Test Cases becomes the active tab, selection A9
$tc.Range("A9").Select()
